$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TPM-derived NATMI metrics (ligand/receptor specificity and edge weights)
# recomputed with new TPM values for the Ctf1-Il6st ligand-receptor pair.
$ws.Cells.Item(2, 9).Value = 0.06742302838872502
$ws.Cells.Item(2, 10).Value = 0.06742302838872503
$ws.Cells.Item(2, 13).Value = 82.43338033333333
$ws.Cells.Item(2, 14).Value = 247.300141
$ws.Cells.Item(2, 15).Value = 0.3670006993429558
$ws.Cells.Item(2, 16).Value = 0.3670006993429557
$ws.Cells.Item(2, 17).Value = 11.96534254435055
$ws.Cells.Item(2, 18).Value = 107.688082899155
$ws.Cells.Item(2, 19).Value = 0.02474429857048204
$ws.Cells.Item(2, 20).Value = 0.02474429857048204
$ws.Cells.Item(3, 9).Value = 0.06742302838872502
$ws.Cells.Item(3, 10).Value = 0.06742302838872503
$ws.Cells.Item(3, 15).Value = 0.3956886215996139
$ws.Cells.Item(3, 16).Value = 0.3956886215996139
$ws.Cells.Item(3, 17).Value = 12.90065633885055
$ws.Cells.Item(3, 19).Value = 0.02667852516720624
$ws.Cells.Item(3, 20).Value = 0.02667852516720624
$ws.Cells.Item(4, 9).Value = 0.06742302838872502
$ws.Cells.Item(4, 10).Value = 0.06742302838872503
$ws.Cells.Item(4, 13).Value = 42.93483766666667
$ws.Cells.Item(4, 14).Value = 128.804513
$ws.Cells.Item(4, 15).Value = 0.1911496942879982
$ws.Cells.Item(4, 16).Value = 0.1911496942879981
$ws.Cells.Item(4, 17).Value = 6.232063245379445
$ws.Cells.Item(4, 18).Value = 56.088569208415
$ws.Cells.Item(4, 19).Value = 0.01288789126447581
$ws.Cells.Item(4, 20).Value = 0.01288789126447581
$ws.Cells.Item(5, 9).Value = 0.06742302838872502
$ws.Cells.Item(5, 10).Value = 0.06742302838872503
$ws.Cells.Item(5, 13).Value = 10.368389
$ws.Cells.Item(5, 14).Value = 31.105167
$ws.Cells.Item(5, 15).Value = 0.04616098476943217
$ws.Cells.Item(5, 16).Value = 0.04616098476943217
$ws.Cells.Item(5, 17).Value = 1.504988943998333
$ws.Cells.Item(5, 18).Value = 13.544900495985
$ws.Cells.Item(5, 19).Value = 0.003112313386560928
$ws.Cells.Item(5, 20).Value = 0.003112313386560929
$ws.Cells.Item(6, 7).Value = 0.3560033333333333
$ws.Cells.Item(6, 9).Value = 0.1653637426357309
$ws.Cells.Item(6, 10).Value = 0.1653637426357309
$ws.Cells.Item(6, 13).Value = 82.43338033333333
$ws.Cells.Item(6, 14).Value = 247.300141
$ws.Cells.Item(6, 15).Value = 0.3670006993429558
$ws.Cells.Item(6, 16).Value = 0.3670006993429557
$ws.Cells.Item(6, 17).Value = 29.34655817660111
$ws.Cells.Item(6, 18).Value = 264.1190235894099
$ws.Cells.Item(6, 19).Value = 0.0606886091932818
$ws.Cells.Item(6, 20).Value = 0.0606886091932818
$ws.Cells.Item(7, 7).Value = 0.3560033333333333
$ws.Cells.Item(7, 9).Value = 0.1653637426357309
$ws.Cells.Item(7, 10).Value = 0.1653637426357309
$ws.Cells.Item(7, 15).Value = 0.3956886215996139
$ws.Cells.Item(7, 16).Value = 0.3956886215996139
$ws.Cells.Item(7, 17).Value = 31.6405368556011
$ws.Cells.Item(7, 19).Value = 0.06543255138608568
$ws.Cells.Item(7, 20).Value = 0.06543255138608568
$ws.Cells.Item(8, 7).Value = 0.3560033333333333
$ws.Cells.Item(8, 9).Value = 0.1653637426357309
$ws.Cells.Item(8, 10).Value = 0.1653637426357309
$ws.Cells.Item(8, 13).Value = 42.93483766666667
$ws.Cells.Item(8, 14).Value = 128.804513
$ws.Cells.Item(8, 15).Value = 0.1911496942879982
$ws.Cells.Item(8, 16).Value = 0.1911496942879981
$ws.Cells.Item(8, 17).Value = 15.28494532545889
$ws.Cells.Item(8, 18).Value = 137.56450792913
$ws.Cells.Item(8, 19).Value = 0.03160922885113918
$ws.Cells.Item(8, 20).Value = 0.03160922885113918
$ws.Cells.Item(9, 7).Value = 0.3560033333333333
$ws.Cells.Item(9, 9).Value = 0.1653637426357309
$ws.Cells.Item(9, 10).Value = 0.1653637426357309
$ws.Cells.Item(9, 13).Value = 10.368389
$ws.Cells.Item(9, 14).Value = 31.105167
$ws.Cells.Item(9, 15).Value = 0.04616098476943217
$ws.Cells.Item(9, 16).Value = 0.04616098476943217
$ws.Cells.Item(9, 17).Value = 3.691181045296666
$ws.Cells.Item(9, 18).Value = 33.22062940767
$ws.Cells.Item(9, 19).Value = 0.007633353205224277
$ws.Cells.Item(9, 20).Value = 0.007633353205224277
$ws.Cells.Item(10, 7).Value = 1.651695
$ws.Cells.Item(10, 8).Value = 4.955085
$ws.Cells.Item(10, 9).Value = 0.767213228975544
$ws.Cells.Item(10, 10).Value = 0.7672132289755441
$ws.Cells.Item(10, 13).Value = 82.43338033333333
$ws.Cells.Item(10, 14).Value = 247.300141
$ws.Cells.Item(10, 15).Value = 0.3670006993429558
$ws.Cells.Item(10, 16).Value = 0.3670006993429557
$ws.Cells.Item(10, 17).Value = 136.154802129665
$ws.Cells.Item(10, 18).Value = 1225.393219166985
$ws.Cells.Item(10, 19).Value = 0.2815677915791919
$ws.Cells.Item(10, 20).Value = 0.2815677915791919
$ws.Cells.Item(11, 7).Value = 1.651695
$ws.Cells.Item(11, 8).Value = 4.955085
$ws.Cells.Item(11, 9).Value = 0.767213228975544
$ws.Cells.Item(11, 10).Value = 0.7672132289755441
$ws.Cells.Item(11, 15).Value = 0.3956886215996139
$ws.Cells.Item(11, 16).Value = 0.3956886215996139
$ws.Cells.Item(11, 17).Value = 146.797829201165
$ws.Cells.Item(11, 18).Value = 1321.180462810485
$ws.Cells.Item(11, 19).Value = 0.303577545046322
$ws.Cells.Item(11, 20).Value = 0.303577545046322
$ws.Cells.Item(12, 7).Value = 1.651695
$ws.Cells.Item(12, 8).Value = 4.955085
$ws.Cells.Item(12, 9).Value = 0.767213228975544
$ws.Cells.Item(12, 10).Value = 0.7672132289755441
$ws.Cells.Item(12, 13).Value = 42.93483766666667
$ws.Cells.Item(12, 14).Value = 128.804513
$ws.Cells.Item(12, 15).Value = 0.1911496942879982
$ws.Cells.Item(12, 16).Value = 0.1911496942879981
$ws.Cells.Item(12, 17).Value = 70.915256699845
$ws.Cells.Item(12, 18).Value = 638.2373102986051
$ws.Cells.Item(12, 19).Value = 0.1466525741723832
$ws.Cells.Item(12, 20).Value = 0.1466525741723832
$ws.Cells.Item(13, 7).Value = 1.651695
$ws.Cells.Item(13, 8).Value = 4.955085
$ws.Cells.Item(13, 9).Value = 0.767213228975544
$ws.Cells.Item(13, 10).Value = 0.7672132289755441
$ws.Cells.Item(13, 13).Value = 10.368389
$ws.Cells.Item(13, 14).Value = 31.105167
$ws.Cells.Item(13, 15).Value = 0.04616098476943217
$ws.Cells.Item(13, 16).Value = 0.04616098476943217
$ws.Cells.Item(13, 17).Value = 17.125416269355
$ws.Cells.Item(13, 18).Value = 154.128746424195
$ws.Cells.Item(13, 19).Value = 0.03541531817764697
$ws.Cells.Item(13, 20).Value = 0.03541531817764697
